# Update cryptos sheet with refreshed prices / 1h volume percentages.
# Rows 29 and 30 also swap coin identity (ImmutableX <-> EthereumClassic).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "38.844.39"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "'" + "2.142.80"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'" + "228.31"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("D7").Value = "'" + "62.21"
$ws.Range("E7").Value = "  +2.58%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("D10").Value = "'" + "0.0846"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "'" + "15.99"
$ws.Range("E12").Value = "  +6.76%  "
$ws.Range("D13").Value = "'" + "2.458.07"
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("D14").Value = "'" + "22.15"
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("D16").Value = "'" + "5.52"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").Value = "'" + "2.140.81"
$ws.Range("E17").Value = "  +2.52%  "
$ws.Range("D18").Value = "'" + "38.869.29"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").Value = "'" + "71.79"
$ws.Range("E20").Value = "  +1.65%  "
$ws.Range("D21").Value = "'" + "0.0₃0846"
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("D22").Value = "'" + "227.62"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'" + "2.37"
$ws.Range("E24").Value = "  -2.69%  "
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").Value = "'" + "9.54"
$ws.Range("E26").Value = "  +1.24%  "
$ws.Range("D27").Value = "'" + "170.51"
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'" + "19.55"
$ws.Range("E29").Value = "  +2.00%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'" + "1.41"
$ws.Range("E30").Value = "  -3.05%  "
$ws.Range("E31").Value = "  +9.61%  "
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("D33").Value = "'" + "4.59"
$ws.Range("E33").Value = "  +2.19%  "
$ws.Range("D34").Value = "'" + "4.81"
$ws.Range("E34").Value = "  +2.30%  "
$ws.Range("D35").Value = "'" + "7.16"
$ws.Range("E35").Value = "  +11.35%  "
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("D37").Value = "'" + "2.41"
$ws.Range("E37").Value = "  +0.65%  "
$ws.Range("D38").Value = "'" + "3.55"
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("D39").Value = "'" + "0.999"
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("D40").Value = "'" + "18.15"
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("E41").Value = "  +3.32%  "
$ws.Range("D42").Value = "'" + "102.64"
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("D43").Value = "'" + "1.534.22"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("E44").Value = "  +6.40%  "
$ws.Range("D45").Value = "'" + "7.83"
$ws.Range("E45").Value = "  +2.44%  "
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("E47").Value = "  +6.31%  "
$ws.Range("D48").Value = "'" + "0.0915"
$ws.Range("E48").Value = "  -1.16%  "
$ws.Range("D49").Value = "'" + "4.15"
$ws.Range("E49").Value = "  +1.46%  "
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").Value = "'" + "2.342.13"
$ws.Range("E51").Value = "  +2.37%  "